$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2197657.5
$ws.Range("J17").Value = 2197657.5
$ws.Range("L17").Value = 6592972.5
$ws.Range("N17").Value = -6593308.5

$ws.Range("H20").Value = 2364
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H28").Value = 290.3125
$ws.Range("I28").Value = 225
$ws.Range("J28").Value = 486.25
$ws.Range("K28").Value = 225
$ws.Range("L28").Value = 486.25
$ws.Range("M28").Value = 260
$ws.Range("N28").Value = -1456.25

$ws.Range("H35").Value = 2364
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H57").Value = 50333.668
$ws.Range("J57").Value = 50333.668
$ws.Range("L57").Value = 151001.004
$ws.Range("N57").Value = -151999.004

$ws.Range("H116").Value = 9685
$ws.Range("I116").Value = 15986.429
$ws.Range("J116").Value = 2333.3333
$ws.Range("K116").Value = 15986.429
$ws.Range("L116").Value = 2333.3333
$ws.Range("M116").Value = -12544.429
$ws.Range("N116").Value = -9217.3333

$ws.Range("H132").Value = 2560.3845
$ws.Range("I132").Value = 2523.75
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 7571.25
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5041.25
$ws.Range("N132").Value = -14060

$ws.Range("H138").Value = 5027.8115
$ws.Range("I138").Value = 948.0345
$ws.Range("J138").Value = 9957.541999999999
$ws.Range("K138").Value = 2844.1035
$ws.Range("L138").Value = 29872.626
$ws.Range("M138").Value = 2295.8965
$ws.Range("N138").Value = -40152.626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1515.6478
$ws.Range("I74").Value = 1503.1936
$ws.Range("K74").Value = 1503.1936
$ws.Range("M74").Value = -629.1936000000001

$ws.Range("H77").Value = 1515.6478
$ws.Range("I77").Value = 1503.1936
$ws.Range("K77").Value = 7515.968000000001
$ws.Range("M77").Value = -3147.968000000001

$ws.Range("H97").Value = 1401.5714
$ws.Range("I97").Value = 950
$ws.Range("J97").Value = 1582.2
$ws.Range("K97").Value = 950
$ws.Range("L97").Value = 1582.2
$ws.Range("M97").Value = -454
$ws.Range("N97").Value = -2574.2

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H122").Value = 3207232.2
$ws.Range("I122").Value = 3207232.2
$ws.Range("K122").Value = 9621696.600000001
$ws.Range("M122").Value = -9619246.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 46740
$ws.Range("J59").Value = 46740
$ws.Range("L59").Value = 46740
$ws.Range("N59").Value = -48434

$ws.Range("H86").Value = 22224062
$ws.Range("I86").Value = 30304692
$ws.Range("J86").Value = 2326.75
$ws.Range("K86").Value = 30304692
$ws.Range("L86").Value = 2326.75
$ws.Range("M86").Value = -30303569
$ws.Range("N86").Value = -4572.75

$ws.Range("H89").Value = 22224062
$ws.Range("I89").Value = 30304692
$ws.Range("J89").Value = 2326.75
$ws.Range("K89").Value = 151523460
$ws.Range("L89").Value = 11633.75
$ws.Range("M89").Value = -151517844
$ws.Range("N89").Value = -22865.75

$ws.Range("H94").Value = 1465.5385
$ws.Range("I94").Value = 927.44446
$ws.Range("K94").Value = 927.44446
$ws.Range("M94").Value = -476.44446

$ws.Range("H102").Value = 32556
$ws.Range("I102").Value = 32556
$ws.Range("K102").Value = 32556
$ws.Range("M102").Value = -29311

$ws.Range("H105").Value = 18131.77
$ws.Range("I105").Value = 34835.668
$ws.Range("J105").Value = 3814.1428
$ws.Range("K105").Value = 34835.668
$ws.Range("L105").Value = 3814.1428
$ws.Range("M105").Value = -33088.668
$ws.Range("N105").Value = -7308.1428

$ws.Range("H127").Value = 55652
$ws.Range("J127").Value = 55652
$ws.Range("L127").Value = 55652
$ws.Range("N127").Value = -65572

$ws.Range("H133").Value = 30757.5
$ws.Range("J133").Value = 30757.5
$ws.Range("L133").Value = 30757.5
$ws.Range("N133").Value = -40877.5

$ws.Range("H134").Value = 5160.8237
$ws.Range("I134").Value = 6524.5454
$ws.Range("J134").Value = 2660.6667
$ws.Range("K134").Value = 19573.6362
$ws.Range("L134").Value = 7982.000100000001
$ws.Range("M134").Value = -17038.6362
$ws.Range("N134").Value = -13052.0001

$ws.Range("H140").Value = 46054.453
$ws.Range("J140").Value = 46054.453
$ws.Range("L140").Value = 46054.453
$ws.Range("N140").Value = -56414.453

$ws.Range("H141").Value = 40390
$ws.Range("J141").Value = 40390
$ws.Range("L141").Value = 40390
$ws.Range("N141").Value = -50750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40368

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0

$ws.Range("H135").Value = 33233.156
$ws.Range("J135").Value = 33233.156
$ws.Range("L135").Value = 33233.156
$ws.Range("N135").Value = -43373.156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3700.75
$ws.Range("J82").Value = 3945.2778
$ws.Range("L82").Value = 11835.8334
$ws.Range("N82").Value = -12647.8334

$ws.Range("H85").Value = 3700.75
$ws.Range("J85").Value = 3945.2778
$ws.Range("L85").Value = 11835.8334
$ws.Range("N85").Value = -14643.8334

$ws.Range("H121").Value = 932.7593000000001
$ws.Range("I121").Value = 408.75
$ws.Range("J121").Value = 1023.8913
$ws.Range("K121").Value = 1226.25
$ws.Range("L121").Value = 3071.6739
$ws.Range("M121").Value = 83.75
$ws.Range("N121").Value = -5691.6739

$ws.Range("H131").Value = 2565065
$ws.Range("I131").Value = 7693115
$ws.Range("J131").Value = 1039.9615
$ws.Range("K131").Value = 23079345
$ws.Range("L131").Value = 3119.8845
$ws.Range("M131").Value = -23074305
$ws.Range("N131").Value = -13199.8845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.77778000000001
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = 99.40000000000001
$ws.Range("K2").Value = 35
$ws.Range("L2").Value = 99.40000000000001
$ws.Range("M2").Value = 78
$ws.Range("N2").Value = -325.4

$ws.Range("H97").Value = 406.07144
$ws.Range("I97").Value = 440.41666
$ws.Range("K97").Value = 440.41666
$ws.Range("M97").Value = 55.58334000000002

$ws.Range("H132").Value = 5589.8887
$ws.Range("I132").Value = 6702.2
$ws.Range("J132").Value = 4199.5
$ws.Range("K132").Value = 20106.6
$ws.Range("L132").Value = 12598.5
$ws.Range("M132").Value = -17576.6
$ws.Range("N132").Value = -17658.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 69080.60000000001
$ws.Range("I7").Value = 113244.89
$ws.Range("K7").Value = 113244.89
$ws.Range("M7").Value = -113132.89

$ws.Range("H16").Value = 1522.3334
$ws.Range("J16").Value = 2100.2
$ws.Range("L16").Value = 2100.2
$ws.Range("N16").Value = -2440.2

$ws.Range("H93").Value = 100040580
$ws.Range("I93").Value = 100450
$ws.Range("K93").Value = 100450
$ws.Range("M93").Value = -99202

$ws.Range("H94").Value = 24000
$ws.Range("J94").Value = 24000
$ws.Range("L94").Value = 24000
$ws.Range("N94").Value = -25352

$ws.Range("H126").Value = 69080.60000000001
$ws.Range("I126").Value = 113244.89
$ws.Range("K126").Value = 339734.67
$ws.Range("M126").Value = -337264.67

$ws.Range("H132").Value = 10322671
$ws.Range("I132").Value = 12750688
$ws.Range("J132").Value = 3599.75
$ws.Range("K132").Value = 38252064
$ws.Range("L132").Value = 10799.25
$ws.Range("M132").Value = -38249534
$ws.Range("N132").Value = -15859.25

$ws.Range("H136").Value = 4739.0835
$ws.Range("I136").Value = 3965.2632
$ws.Range("J136").Value = 7679.6
$ws.Range("K136").Value = 11895.7896
$ws.Range("L136").Value = 23038.8
$ws.Range("M136").Value = -9345.7896
$ws.Range("N136").Value = -28138.8
